# Update the active workbook from the Set12 "moderate_noisy" dual-NLM
# experiment dump to the "pg_noisy" dual-NLM experiment dump:
#  - clear out the old per-image PSNR/SSIM/Score rows (A2:D12)
#  - leave the header row and the formatted blank spacer row (B16) intact
#  - append the two new summary strings for the pg_noisy run

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old per-image result rows, keeping the header row (row 1)
# and the blank formatted row (row 16) untouched.
$ws.Range("A2:D12").Clear() | Out-Null

# Add the two new summary rows for the pg_noisy dual experiment.
$ws.Range("D20").Value = "PSNR = 22.08 | SSIM = 0.4068 | Score = 31.38"
$ws.Range("D21").Value = "PSNR = 22.78 | SSIM = 0.4333 | Score = 33.05"

# Move the active selection, matching the saved workbook state.
$ws.Range("P7").Select() | Out-Null
